# "them history va options" - clear the stray option marker left in I3
# (shared-string ";1") on Sheet1. Removing the cell also drops the now-
# unused shared-string entry so sharedStrings.xml count/uniqueCount shrink
# from 5 to 4, matching the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").ClearContents()
